# Applies the "Diverse kalenderarbeider og kontingent" edit:
#  - adds two new shared strings: "frmInnsatsScr" and "frmDugnadInnsatsScr"
#  - writes them into column C of the flag-description table (rows 23 and 13)
#  - updates the active selection / scroll position of the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 (A23 = 22) gets the first newly introduced label -> shared string index 11
$ws.Range("C23").Value2 = "frmInnsatsScr"

# Row 13 (A13 = 12) gets the second newly introduced label -> shared string index 12
$ws.Range("C13").Value2 = "frmDugnadInnsatsScr"

# Update the view: scrolled down a bit, with E19 as the active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("E19").Select()
